# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets,
# reflecting newly generated output data (commit: "Update gh-pages to output
# generated at a3196b5").

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 289
    $ws.Range("F5").Value = 267
}
